$wb = $excel.ActiveWorkbook

# --- Update the conversion text on "Hoja1" (sheet1), cell A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.03 = 7689.86 pesos`n✅ 7689.86 pesos = 2.02 = 918.1 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on "tasas" (sheet2) ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 493
$ws2.Range("O10").Value = 3791.1
$ws2.Range("N12").Value = 3811
$ws2.Range("O12").Value = 455
